$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $range = $d.Content
    $range.Find.Execute($findText, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $replaceText, 2)
}

# 1. Update activation date
Replace-Text "Ativação: 01/01/2018" "Ativação: 01/01/2025"

# 2. Prepend new sentences to the "Programa resumido" paragraph
Replace-Text "A disciplina aborda uma série de processos industriais que são utilizados no fracionamento e na conversão da biomassa vegetal ligninficada em produtos elaborados como celulose e papel, derivados de celulose, carvão e açúcares." `
             "Estrutura e ultraestrutura dos materiais lignocelulósicos, celulose, hemiceluloses e outras polioses. Lignina, extrativos e composição da casca. Reações em meio ácido, meio alcalino. A disciplina aborda uma série de processos industriais que são utilizados no fracionamento e na conversão da biomassa vegetal ligninficada em produtos elaborados como celulose e papel, derivados de celulose, carvão e açúcares."

# 3. Insert spaces between numbered items in the "Programa" paragraph (Portuguese version only)
Replace-Text "química.2. Produção" "química. 2. Produção"
Replace-Text "de papel.3. Produção" "de papel. 3. Produção"
Replace-Text "epóxidos.4. Conversão" "epóxidos. 4. Conversão"
Replace-Text "vegetal.5. Produção" "vegetal. 5. Produção"
Replace-Text "celular.6. Processos" "celular. 6. Processos"

# 4. Update the NF formula with the Estudo de Caso weighting
Replace-Text "NF=(P1+P2)/2" "NF=(P1+P2)/2 x 0,9 + Estudo de Caso x 0,1."

# 5. Add trailing period to the MR formula
Replace-Text "MR=(NF=PR)/2" "MR=(NF=PR)/2."

$word.ActiveDocument.Save()
